$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-08-08 Friday"

# Update the division problems in the table, addressed by (row, column)
# to avoid any ambiguity from values that coincide with other cells old/new text.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "36÷2=18, 0"
$t.Cell(1, 2).Range.Text = "66÷9=7, 3"
$t.Cell(1, 3).Range.Text = "40÷4=10, 0"
$t.Cell(1, 4).Range.Text = "79÷8=9, 7"
$t.Cell(1, 5).Range.Text = "20÷6=3, 2"
$t.Cell(5, 1).Range.Text = "84÷6=14, 0"
$t.Cell(5, 2).Range.Text = "85÷5=17, 0"
$t.Cell(5, 3).Range.Text = "59÷2=29, 1"
$t.Cell(5, 4).Range.Text = "26÷3=8, 2"
$t.Cell(5, 5).Range.Text = "17÷2=8, 1"
$t.Cell(9, 1).Range.Text = "14÷6=2, 2"
$t.Cell(9, 2).Range.Text = "81÷9=9, 0"
$t.Cell(9, 3).Range.Text = "87÷9=9, 6"
$t.Cell(9, 4).Range.Text = "26÷3=8, 2"
$t.Cell(9, 5).Range.Text = "30÷8=3, 6"
$t.Cell(13, 1).Range.Text = "39÷2=19, 1"
$t.Cell(13, 2).Range.Text = "95÷2=47, 1"
$t.Cell(13, 3).Range.Text = "75÷6=12, 3"
$t.Cell(13, 4).Range.Text = "82÷4=20, 2"
$t.Cell(13, 5).Range.Text = "31÷4=7, 3"
$t.Cell(17, 1).Range.Text = "96÷6=16, 0"
$t.Cell(17, 2).Range.Text = "65÷2=32, 1"
$t.Cell(17, 3).Range.Text = "79÷3=26, 1"
$t.Cell(17, 4).Range.Text = "27÷9=3, 0"
$t.Cell(17, 5).Range.Text = "38÷3=12, 2"
